$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 8-36: A,B,D,E,F,G,H,Q,R (col C stays "Ovaliderad" for all, unchanged)
$data = @(
    @{Row=8; A=111756142; B=90087; D='LC'; E=3298; F='Trådticka'; G='Climacocystis borealis'; H='(Fr.) Kotl. & Pouzar'; Q=454002.5104495964; R=7073638.391199326},
    @{Row=9; A=111756159; B=89423; D='NT'; E=5432; F='Granticka'; G='Porodaedalea chrysoloma'; H='(Fr.) Fiasson & Niemelä'; Q=453621.4908246148; R=7073983.576241231},
    @{Row=10; A=111756163; B=77515; D='NT'; E=6425; F='Garnlav'; G='Alectoria sarmentosa'; H='(Ach.) Ach.'; Q=453955.6479769219; R=7073945.9492877},
    @{Row=11; A=111756172; B=85715; D='NT'; E=510; F='Doftskinn'; G='Cystostereum murrayi'; H='(Berk. & M.A. Curtis.) Pouzar'; Q=453938.5789576455; R=7073959.46382203},
    @{Row=12; A=111756143; B=90087; D='LC'; E=3298; F='Trådticka'; G='Climacocystis borealis'; H='(Fr.) Kotl. & Pouzar'; Q=453950.9091414157; R=7073591.829928016},
    @{Row=13; A=111756140; B=89405; D='NT'; E=1202; F='Ullticka'; G='Phellinidium ferrugineofuscum'; H='(P.Karst.) Fiasson & Niemelä'; Q=453820.6239011836; R=7074037.242731699},
    @{Row=14; A=111756169; B=77515; D='NT'; E=6425; F='Garnlav'; G='Alectoria sarmentosa'; H='(Ach.) Ach.'; Q=453910.2023238647; R=7073654.334338664},
    @{Row=15; A=111756139; B=89405; D='NT'; E=1202; F='Ullticka'; G='Phellinidium ferrugineofuscum'; H='(P.Karst.) Fiasson & Niemelä'; Q=453692.6056797595; R=7074032.491935454},
    @{Row=16; A=111756170; B=96265; D='LC'; E=219790; F='Fläcknycklar'; G='Dactylorhiza maculata'; H='(L.) Soó'; Q=453738.5427278728; R=7073724.066700204},
    @{Row=17; A=111756147; B=89425; D='NT'; E=5442; F='Tallticka'; G='Porodaedalea pini'; H='(Brot.) Murrill'; Q=453989.3915585176; R=7073710.21875874},
    @{Row=18; A=111756155; B=89423; D='NT'; E=5432; F='Granticka'; G='Porodaedalea chrysoloma'; H='(Fr.) Fiasson & Niemelä'; Q=453863.4009631127; R=7073965.428905412},
    @{Row=19; A=111756160; B=77515; D='NT'; E=6425; F='Garnlav'; G='Alectoria sarmentosa'; H='(Ach.) Ach.'; Q=453815.5156181521; R=7073870.182023689},
    @{Row=20; A=111756171; B=88899; D='NT'; E=3286; F='Flattoppad klubbsvamp'; G='Clavariadelphus truncatus'; H='(Quél.) Donk'; Q=453750.6060291855; R=7073942.323881648},
    @{Row=21; A=111756157; B=89423; D='NT'; E=5432; F='Granticka'; G='Porodaedalea chrysoloma'; H='(Fr.) Fiasson & Niemelä'; Q=453981.5111392652; R=7073807.172376178},
    @{Row=22; A=111756154; B=96674; D='LC'; E=219880; F='Kransrams'; G='Polygonatum verticillatum'; H='(L.) All.'; Q=453614.9183513908; R=7074108.35826167},
    @{Row=23; A=111756161; B=77515; D='NT'; E=6425; F='Garnlav'; G='Alectoria sarmentosa'; H='(Ach.) Ach.'; Q=453723.2573215028; R=7074069.623294062},
    @{Row=24; A=111756167; B=77515; D='NT'; E=6425; F='Garnlav'; G='Alectoria sarmentosa'; H='(Ach.) Ach.'; Q=454002.5104495964; R=7073638.391199326},
    @{Row=25; A=111756156; B=89423; D='NT'; E=5432; F='Granticka'; G='Porodaedalea chrysoloma'; H='(Fr.) Fiasson & Niemelä'; Q=453978.4965374623; R=7073812.964766338},
    @{Row=26; A=111756151; B=95532; D='LC'; E=221945; F='Revlummer'; G='Lycopodium annotinum'; H='L.'; Q=453609.4901279925; R=7074130.545069677},
    @{Row=27; A=111756166; B=77515; D='NT'; E=6425; F='Garnlav'; G='Alectoria sarmentosa'; H='(Ach.) Ach.'; Q=453981.6720900657; R=7073697.065866594},
    @{Row=28; A=111756153; B=96674; D='LC'; E=219880; F='Kransrams'; G='Polygonatum verticillatum'; H='(L.) All.'; Q=453707.5163784204; R=7073721.869806641},
    @{Row=29; A=111756162; B=77515; D='NT'; E=6425; F='Garnlav'; G='Alectoria sarmentosa'; H='(Ach.) Ach.'; Q=453922.6243923472; R=7073958.370937477},
    @{Row=30; A=111756165; B=77515; D='NT'; E=6425; F='Garnlav'; G='Alectoria sarmentosa'; H='(Ach.) Ach.'; Q=453984.2379404157; R=7073751.417626478},
    @{Row=31; A=111756158; B=89423; D='NT'; E=5432; F='Granticka'; G='Porodaedalea chrysoloma'; H='(Fr.) Fiasson & Niemelä'; Q=454002.8592168373; R=7073783.424762985},
    @{Row=32; A=111756164; B=77515; D='NT'; E=6425; F='Garnlav'; G='Alectoria sarmentosa'; H='(Ach.) Ach.'; Q=453971.0747186596; R=7073820.148138274},
    @{Row=33; A=111756150; B=95532; D='LC'; E=221945; F='Revlummer'; G='Lycopodium annotinum'; H='L.'; Q=453976.2702886119; R=7073812.112971266},
    @{Row=34; A=111756148; B=96266; D='LC'; E=223591; F='Skogsnycklar'; G='Dactylorhiza maculata subsp. fuchsii'; H='(Druce) Hyl.'; Q=453747.0542679164; R=7073851.289854143},
    @{Row=35; A=111756141; B=89405; D='NT'; E=1202; F='Ullticka'; G='Phellinidium ferrugineofuscum'; H='(P.Karst.) Fiasson & Niemelä'; Q=453610.1793069927; R=7074087.205471905},
    @{Row=36; A=111756168; B=77515; D='NT'; E=6425; F='Garnlav'; G='Alectoria sarmentosa'; H='(Ach.) Ach.'; Q=453958.9423245317; R=7073596.134472342}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
    $ws.Cells.Item($r, 7).Value = $item.G
    $ws.Cells.Item($r, 8).Value = $item.H
    $ws.Cells.Item($r, 17).Value = $item.Q
    $ws.Cells.Item($r, 18).Value = $item.R
}